$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.190.37"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.84%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.824.98"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.79%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4232"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.93%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3689"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07237"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.52%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8545"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.96"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.87%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.828.90"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.699"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.07071"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.78%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.303"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.74%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "89.84"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.32%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.10%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008848"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.97%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "27.244.07"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.115"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.39%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.89"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.25%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.052.09"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.978"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.56%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "152.37"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.198"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.39%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.42"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.241"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.60%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "116.44"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.17%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08839"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.14%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.191"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.81%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7498"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.933"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.440"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.87%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.104"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("E38").Value = "  -0.50%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05236"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.95%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "7.274"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.869"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1698"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.00%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5037"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.33%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.654"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("E45").Value = "  -0.73%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "106.35"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4737"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.98%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06390"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.666"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.877"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "